# Update Nalco aluminium-ingot price sheet with the latest daily record.
#
# The sheet is a rolling daily log (newest date at row 2, oldest at the
# bottom, row 1 is the header). This commit adds a new "today" row right
# under the header, which pushes every existing data row down by one; the
# row that used to be last is *also* preserved by duplicating it onto a
# new row at the bottom (the historical window grows by one row instead
# of the oldest entry being dropped).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Capture the data we need before any rows move ---------------------

# Current newest row (row 2) - becomes the template for the row we insert,
# except the Date column (A) advances by one day.
$oldA2 = $ws.Range("A2").Value2
$oldB2 = $ws.Range("B2").Value2
$oldC2 = $ws.Range("C2").Value2
$oldD2 = $ws.Range("D2").Value2
$oldE2 = $ws.Range("E2").Value2
$oldF2 = $ws.Range("F2").Value2

# Current last row - gets duplicated onto the new bottom row, unchanged.
$lastRow = $ws.UsedRange.Rows.Count
$oldALast = $ws.Cells.Item($lastRow, 1).Value2
$oldBLast = $ws.Cells.Item($lastRow, 2).Value2
$oldCLast = $ws.Cells.Item($lastRow, 3).Value2
$oldDLast = $ws.Cells.Item($lastRow, 4).Value2
$oldELast = $ws.Cells.Item($lastRow, 5).Value2
$oldFLast = $ws.Cells.Item($lastRow, 6).Value2

$newDate = [DateTime]::ParseExact($oldA2, "dd-MM-yyyy", $null).AddDays(1).ToString("dd-MM-yyyy")

$newLastRow = $lastRow + 1

# --- Insert the new top row, pushing everything else down --------------

$ws.Rows("2:2").Insert()

# Write the new row's values as plain text where applicable (forcing the
# Text number format first prevents Excel from "helpfully" reinterpreting
# dd-MM-yyyy strings as dates).
$ws.Range("A2").NumberFormat = "@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("F2").NumberFormat = "@"

$ws.Range("A2").Value = $newDate
$ws.Range("B2").Value = $oldB2
$ws.Range("C2").Value = $oldC2
$ws.Range("D2").Value = $oldD2
$ws.Range("E2").Value = $oldE2
$ws.Range("F2").Value = $oldF2

$ws.Hyperlinks.Add($ws.Range("F2"), $oldF2) | Out-Null

# Re-apply the same formatting (alignment, 0.000 numeric format, etc.) as
# the rest of the table by copying it from the row right below, which is
# untouched data and still carries the original styling. Doing this after
# Hyperlinks.Add also clears the blue/underline "Hyperlink" style it
# applies, matching how the source sheet formats its link cells as plain
# text.
$ws.Range("A3:F3").Copy()
$ws.Range("A2:F2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Append the duplicated last row at the new bottom -------------------

$ws.Cells.Item($newLastRow, 1).NumberFormat = "@"
$ws.Cells.Item($newLastRow, 2).NumberFormat = "@"
$ws.Cells.Item($newLastRow, 3).NumberFormat = "@"
$ws.Cells.Item($newLastRow, 5).NumberFormat = "@"
$ws.Cells.Item($newLastRow, 6).NumberFormat = "@"

$ws.Cells.Item($newLastRow, 1).Value = $oldALast
$ws.Cells.Item($newLastRow, 2).Value = $oldBLast
$ws.Cells.Item($newLastRow, 3).Value = $oldCLast
$ws.Cells.Item($newLastRow, 4).Value = $oldDLast
$ws.Cells.Item($newLastRow, 5).Value = $oldELast
$ws.Cells.Item($newLastRow, 6).Value = $oldFLast

$ws.Hyperlinks.Add($ws.Cells.Item($newLastRow, 6), $oldFLast) | Out-Null

$ws.Range("A" + ($newLastRow - 1) + ":F" + ($newLastRow - 1)).Copy()
$ws.Range("A" + $newLastRow + ":F" + $newLastRow).PasteSpecial(-4122)
$excel.CutCopyMode = $false
